$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style (bold, thin border, centered) from an existing header cell (A1)
# onto the three new header cells before setting their text.
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats

# New header labels
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in the team record (Wins/Losses/Ties) for every data row (2-41)
$lastRow = 41
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 96   # AD - Wins
    $ws.Cells.Item($r, 31).Value = 66   # AE - Losses
    $ws.Cells.Item($r, 32).Value = 0    # AF - Ties
}
